$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 151, shifting the existing rows 151-196 down to 152-197
$ws.Rows.Item(151).Insert()

# Populate the newly inserted row 151 with the new record
$ws.Range("A151").Value = 9
$ws.Range("B151").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C151").Value = "Metropolitana"
$ws.Range("D151").Value = 44551
$ws.Range("E151").Value = 13
$ws.Range("F151").Value = 100112026
$ws.Range("G151").Value = "Haba"
$ws.Range("H151").Value = "Sin especificar"
$ws.Range("I151").Value = "Primera"
$ws.Range("J151").Value = 61
$ws.Range("K151").Value = 16000
$ws.Range("L151").Value = 18000
$ws.Range("M151").Value = 17016
$ws.Range("N151").Value = "`$/saco 25 kilos"
$ws.Range("O151").Value = "Carahue"
$ws.Range("P151").Value = 681
$ws.Range("Q151").Value = 25
$ws.Range("R151").Value = "Hortaliza"
